$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.695.24"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +1.55%  "

# Row 3: 'Ethereum'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.880.15"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  +2.23%  "

# Row 4: 'TetherUSD'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5: 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.24"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.95%  "

# Row 6: 'USDC'
$ws.Range("E6").Value = "  +0.09%  "

# Row 7: 'XRP'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4760"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +2.24%  "

# Row 8: 'Cardano'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2835"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +4.28%  "

# Row 9: 'Dogecoin'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06493"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +3.57%  "

# Row 10: 'Solana'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.60"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +15.51%  "

# Row 11: 'WrappedEther'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.877.95"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +2.38%  "

# Row 12: 'TRON'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07563"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  +1.91%  "

# Row 13: 'Litecoin'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.27"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +14.09%  "

# Row 14: 'Polkadot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.049"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +2.37%  "

# Row 15: 'Polygon'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6487"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  +4.85%  "

# Row 16: 'BitcoinCash'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "297.74"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +31.93%  "

# Row 17: 'WrappedBTC'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.690.91"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +1.76%  "

# Row 18: 'Avalanche' -> 'Dai'
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +0.17%  "

# Row 19: 'Dai' -> 'Avalanche'
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.07"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +6.23%  "

# Row 20: 'ShibaInu'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007468"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +2.72%  "

# Row 21: 'WrappedliquidstakedEther2.0'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.147.05"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +3.68%  "

# Row 22: 'BinanceUSD'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +0.11%  "

# Row 23: 'Uniswap'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.123"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +5.03%  "

# Row 24: 'Chainlink'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.110"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +4.66%  "

# Row 25: 'Monero'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.46"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +3.22%  "

# Row 26: 'Cosmos'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.214"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +0.34%  "

# Row 27: 'EthereumClassic'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.49"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +9.76%  "

# Row 28: 'LidoDAOToken'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.948"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +4.76%  "

# Row 29: 'Stellar'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1056"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +2.32%  "

# Row 30: 'Toncoin'
$ws.Range("E30").Value = "  -1.33%  "

# Row 31: 'InternetComputer(DFINITY)'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.147"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +2.05%  "

# Row 32: 'Filecoin'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.943"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +3.86%  "

# Row 33: 'Hedera'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04982"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +3.43%  "

# Row 34: 'ARBITRUM'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.168"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +2.69%  "

# Row 35: 'ImmutableX'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7173"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +2.14%  "

# Row 36: 'HuobiToken'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.721"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +0.88%  "

# Row 37: 'VeChain'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01907"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +2.02%  "

# Row 38: 'MXToken'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.702"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +1.99%  "

# Row 39: 'RenderToken'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.046"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +6.92%  "

# Row 40: 'TrustWalletToken'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8944"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +0.53%  "

# Row 41: 'Quant'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.81"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +2.21%  "

# Row 42: 'PaxDollar'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -0.06%  "

# Row 43: 'TheSandbox'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4172"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +4.17%  "

# Row 44: 'FraxShare'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.561"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +0.76%  "

# Row 45: 'Aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.44"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +9.41%  "

# Row 46: 'Aptos'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.312"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +4.40%  "

# Row 47: 'Algorand'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1211"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +1.45%  "

# Row 48: 'EnergySwap' -> 'Elrond'
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.47"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  +4.75%  "

# Row 49: 'Elrond' -> 'EnergySwap'
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.798"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +2.83%  "

# Row 50: 'Cronos'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05611"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +1.88%  "

# Row 51: 'NEARProtocol'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.380"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +2.35%  "
